$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Replace the text of the italic "Read our review..." teaser
#         paragraph (near the end of the document) with the new
#         image-generation prompt text. We locate it with Find (no
#         replacement text supplied) and then assign the new text
#         directly to the found Range so that straight apostrophes
#         are preserved verbatim (Find's own Replacement text field
#         smart-quotes apostrophes, which we must avoid) and the
#         existing run formatting (italic) carries over automatically.
# ------------------------------------------------------------------
$oldTeaser = "Read our review of ARRR! 10K Ways and play for free. Discover its 10,000 ways to win, low volatility, and multiple bonus features."
$newImagePrompt = "Create a feature image for ARRR! 10K Ways that captures the adventurous spirit of a pirate-themed game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing on a sandy beach with palm trees in the background, holding a treasure chest filled with gold coins and jewels. The warrior's clothing and accessories should suggest that they are a pirate on their quest for treasure. The image should also include the game title, ARRR! 10K Ways, in bold and eye-catching letters. The overall design should be bright, colorful, and playful to attract players' attention and generate excitement about the game."

$teaserRange = $d.Content
$found = $teaserRange.Find.Execute($oldTeaser)
$teaserRange.Text = $newImagePrompt

# ------------------------------------------------------------------
# Step 2: Remove the now-duplicate bold "Play ARRR! 10K Ways for
#         Free - Review and Bonus Features" paragraph that used to
#         sit right above the teaser (it is being replaced by the new
#         "Meta description" paragraph inserted under the main title).
# ------------------------------------------------------------------
$needle = "Play ARRR! 10K Ways for Free - Review and Bonus Features"
$dupIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd()
    if ($txt -eq $needle) {
        $dupIndex = $i
    }
}
if ($dupIndex -ge 2) {
    $d.Paragraphs($dupIndex).Range.Delete()
}

# ------------------------------------------------------------------
# Step 3: Insert a new paragraph right after the main Heading1 title
#         containing: an empty leading run, a bold "Meta description"
#         run, and a plain run with the rest of the meta description.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"
$metaPara.Range.InsertAfter("Meta description")
$d.Paragraphs(2).Range.InsertAfter(": Read our review of ARRR! 10K Ways and play for free. Discover its 10,000 ways to win, low volatility, and multiple bonus features.")

# Bold just the "Meta description" label.
$labelStart = $d.Paragraphs(2).Range.Start
$labelLen = "Meta description".Length
$labelRange = $d.Range($labelStart, $labelStart + $labelLen)
$labelRange.Font.Bold = $true

# Split off a genuinely empty leading run (matches the leading empty
# <w:r/> pattern used throughout the rest of the document's body
# paragraphs).
$zeroRange = $d.Range($labelStart, $labelStart)
$zeroRange.InsertBefore("")
